$d = $word.ActiveDocument

# Replace the ID placeholder text (this also consumes the following
# run's trailing space, collapsing the two runs into one).
$d.Content.Find.Execute("**ID__AFFARS_5303_topic_16__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5303_705__ID**", 2)

# Update the first paragraph's formatting: add a paragraph border
# (space=5 on all sides) and change the left indent from 120 to 225 twips
# (i.e. 6pt -> 11.25pt).
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.LeftIndent = 11.25
$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5
